# Regenerate save_data to use K (strikeouts) instead of Strike# (pitch count
# for strikes) in column G. Values were recalculated and are written below
# for each data row (rows 2-34 correspond to the 33 game log entries).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 10
    3  = 7
    4  = 8
    5  = 11
    6  = 10
    7  = 2
    8  = 10
    9  = 11
    10 = 7
    11 = 3
    12 = 10
    13 = 10
    14 = 7
    15 = 7
    16 = 12
    17 = 7
    18 = 4
    19 = 10
    20 = 8
    21 = 8
    22 = 5
    23 = 6
    24 = 2
    25 = 6
    26 = 3
    27 = 12
    28 = 3
    29 = 9
    30 = 8
    31 = 9
    32 = 12
    33 = 4
    34 = 4
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
